$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Mon Oct 09 22:49:06 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 09 22:49:19 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 09 22:49:33 EDT 2023"
